$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Clone row 6 (the "Extension.valueCodeableConcept" slice row) down into a
#    new row 7, preserving both formatting (style s="2", borders, wrap) and
#    the literal values/types (so text cells that look numeric, e.g. the "0"
#    / "1" Min/Max columns, stay shared-string text instead of becoming
#    numbers).
# ---------------------------------------------------------------------------
$ws.Range("A6:AJ6").Copy()
$ws.Range("A7:AJ7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A6:AJ6").Copy()
$ws.Range("A7:AJ7").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

# Row 7 is a hidden detail row, same as the rows above it.
$ws.Rows.Item(7).Hidden = $true

# ---------------------------------------------------------------------------
# 2) Row 7 now equals the OLD row 6 verbatim, except the Path column (A)
#    must read the generic "Extension.value[x]" path instead of the old
#    type-specific path.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Extension.value[x]"

# ---------------------------------------------------------------------------
# 3) Re-purpose row 6 itself into the generic "Extension.value[x]" slicing
#    header row (closed slicing on type, pointing at the value[x] element).
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Extension.value[x]"
$ws.Range("B6").Value = ""
$ws.Range("K6").Value = "Value of extension"
$ws.Range("L6").Value = "Value of extension - may be a resource or one of a constrained set of the data types (see Extensibility in the spec for list)."
$ws.Range("W6").Value = ""
$ws.Range("X6").Value = ""
$ws.Range("Y6").Value = ""
$ws.Range("AA6").Value = "type:`$this}" + "`n"
$ws.Range("AD6").Value = "closed"

# ---------------------------------------------------------------------------
# 4) Sheet-level ranges that now need to span through row 7 instead of 6.
# ---------------------------------------------------------------------------

# Rebuild the AutoFilter over the extended range, restoring both existing
# filter columns (G <> " " / AA blank).
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ7").AutoFilter()
$ws.Range("A1:AJ7").AutoFilter(7, "<> ")
$ws.Range("A1:AJ7").AutoFilter(27, @(""), 7)

# _xlnm._FilterDatabase defined name tracks the same range.
$fdb = $wb.Names.Item(1)
$fdb.RefersTo = "=Elements!`$A`$1:`$AJ`$7"

# Conditional formatting grows by one row (A2:AI5 -> A2:AI6).
$fc = $ws.Range("A2:AI5").FormatConditions().Item(1)
$fc.ModifyAppliesToRange($ws.Range("A2:AI6"))

# Column A ("Path") no longer needs to fit the old 30-character string, so
# its best-fit width shrinks.
$ws.Columns.Item(1).ColumnWidth = 18.1
